# add Save column in s_vals sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", matching style/format of existing header cell G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell H2 = 0 (numeric)
$ws.Range("H2").Value = 0
